{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the target paragraph: \"Today with the indispensable help of Swen ...\"\n// (the final paragraph of the 07/03/2018 dev-diary entry).\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"Today with the indispensable\") === 0) {\n    target = p;\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not locate target paragraph\");\n}\n\n// Replace that single paragraph with the split/annotated version of itself\n// plus the six new diary paragraphs that follow it (including the relocated\n// _GoBack bookmark), reproducing the authored edit in one shot via a\n// flat-OPC OOXML payload (required by insertOoxml).\nconst innerXml =\n  '<w:p><w:r><w:t xml:space=\"preserve\">Today with the indispensable help of Swen the bot now builds at the choke correctly, though the plan needs editing and bit of function need fleshing out, the bot now does more of what its supposed to do.  I was correct in it being a build site/position </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>issue,.</w:t></w:r><w:proofErr w:type=\"gramEnd\"/></w:p>' +\n  '<w:p/>' +\n  '<w:p><w:r><w:t>Altered scouting to happen at start</w:t></w:r><w:r><w:t xml:space=\"preserve\"> of the game.</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\">Changed build </w:t></w:r><w:r><w:t>location</w:t></w:r><w:r><w:t xml:space=\"preserve\"> to change </w:t></w:r><w:r><w:t>correctly</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>Added more functions to scouting so scout only happens once</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\">Added two while loops into the scouting and building positioning so the game has to wait for the probe to arrive at the </w:t></w:r><w:r><w:t>destination</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t xml:space=\"preserve\"> before attempting to fire the next command.</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>Altered the move command in the position function so probe now moves to its destination correctly.</w:t></w:r></w:p>';\n\nconst flatOpcXml =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + innerXml + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ntarget.getRange().insertOoxml(flatOpcXml, \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the target paragraph: \"Today with the indispensable help of Swen ...\"\n# (it is the final paragraph in the dev-diary entry for 07/03/2018).\n$target = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"Today with the indispensable*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not locate target paragraph\"\n}\n\n# Replace that single paragraph with the split/annotated version of itself\n# plus the six new diary paragraphs that follow it (including the relocated\n# _GoBack bookmark), reproducing the authored edit in one shot.\n$xml = @'\n<w:p><w:r><w:t xml:space=\"preserve\">Today with the indispensable help of Swen the bot now builds at the choke correctly, though the plan needs editing and bit of function need fleshing out, the bot now does more of what its supposed to do.  I was correct in it being a build site/position </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>issue,.</w:t></w:r><w:proofErr w:type=\"gramEnd\"/></w:p><w:p/><w:p><w:r><w:t>Altered scouting to happen at start</w:t></w:r><w:r><w:t xml:space=\"preserve\"> of the game.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">Changed build </w:t></w:r><w:r><w:t>location</w:t></w:r><w:r><w:t xml:space=\"preserve\"> to change </w:t></w:r><w:r><w:t>correctly</w:t></w:r></w:p><w:p><w:r><w:t>Added more functions to scouting so scout only happens once</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">Added two while loops into the scouting and building positioning so the game has to wait for the probe to arrive at the </w:t></w:r><w:r><w:t>destination</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t xml:space=\"preserve\"> before attempting to fire the next command.</w:t></w:r></w:p><w:p><w:r><w:t>Altered the move command in the position function so probe now moves to its destination correctly.</w:t></w:r></w:p>\n'@\n\n$target.Range.InsertXML($xml)\n"}
